# New crime data collected — weekly CompStat report update (94th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 31   Number  3" -> "Volume 31   Number  4"
$ws.Range("A8").Characters(21, 1).Text = "4"

# "Report Covering the Week  1/15/2024  Through  1/21/2024"
#   -> "...1/22/2024  Through  1/28/2024"
$ws.Range("C9").Characters(27, 9).Text = "1/22/2024"
$ws.Range("C9").Characters(47, 9).Text = "1/28/2024"

# --- Percent-format helper values ----------------------------------------
$fmtInt = "#,##0"
$fmtPct = "#,##0.0;""-""#,##0.0"

# --- Row 14 (Murder) — only N14 goes from "not applicable" to a number ---
$ws.Range("N14").NumberFormat = $fmtPct
$ws.Range("N14").Value = -100

# --- Row 16 (Rape) --------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -15.384615384615
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = -15.384615384615
$ws.Range("L16").Value = -26.666666666666
$ws.Range("M16").Value = -31.25
$ws.Range("N16").Value = -76.595744680851

# --- Row 17 (Robbery) -----------------------------------------------------
$ws.Range("C17").NumberFormat = $fmtInt
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 11.111111111111
$ws.Range("L17").Value = 11.111111111111
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -47.368421052631

# --- Row 18 (Fel. Assault) -------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = -41.666666666666
$ws.Range("L18").Value = -44
$ws.Range("M18").Value = -12.5
$ws.Range("N18").Value = -87.272727272727

# --- Row 19 (Burglary) ------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -21.428571428571
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 6.976744186046
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 43
$ws.Range("K19").Value = 6.976744186046
$ws.Range("L19").Value = 48.387096774193
$ws.Range("M19").Value = 119.047619047619
$ws.Range("N19").Value = 53.333333333333

# --- Row 20 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -94.117647058823

# --- Row 21 (TOTAL, bold) ----------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -25.925925925925
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = -12.371134020618
$ws.Range("I21").Value = 85
$ws.Range("J21").Value = 97
$ws.Range("K21").Value = -12.371134020618
$ws.Range("L21").Value = -1.162790697674
$ws.Range("M21").Value = 34.920634920634
$ws.Range("N21").Value = -69.202898550724

# --- Row 22 (G.L.A.) -----------------------------------------------------------
$ws.Range("C22").NumberFormat = $fmtInt
$ws.Range("C22").Value = 2
$ws.Range("D22").NumberFormat = $fmtInt
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = $fmtPct
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 50
# L22 and N22 remain "not applicable" (unchanged)
$ws.Range("M22").Value = 200

# --- Row 23 (TOTAL, Transit) -----------------------------------------------------
$ws.Range("C23").NumberFormat = $fmtInt
$ws.Range("C23").Value = 1
$ws.Range("D23").NumberFormat = $fmtInt
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = $fmtPct
$ws.Range("E23").Value = 0
$ws.Range("F23").NumberFormat = $fmtInt
$ws.Range("F23").Value = 1
$ws.Range("G23").NumberFormat = $fmtInt
$ws.Range("G23").Value = 1
$ws.Range("H23").NumberFormat = $fmtPct
$ws.Range("H23").Value = 0
$ws.Range("I23").NumberFormat = $fmtInt
$ws.Range("I23").Value = 1
$ws.Range("J23").NumberFormat = $fmtInt
$ws.Range("J23").Value = 1
$ws.Range("K23").NumberFormat = $fmtPct
$ws.Range("K23").Value = 0
# L23, M23, N23 remain "not applicable" (unchanged)

# --- Row 24 (TOTAL, Housing) -----------------------------------------------------
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 43.75
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 14.035087719298
$ws.Range("I24").Value = 65
$ws.Range("J24").Value = 57
$ws.Range("K24").Value = 14.035087719298
$ws.Range("L24").Value = -12.162162162162
$ws.Range("M24").Value = 80.555555555555
# N24 remains "not applicable" (unchanged)

# --- Row 25 (Petit Larceny) --------------------------------------------------------
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 9.523809523809
$ws.Range("I25").Value = 23
$ws.Range("J25").Value = 21
$ws.Range("K25").Value = 9.523809523809
$ws.Range("L25").Value = 35.294117647058
$ws.Range("M25").Value = 64.285714285714
# N25 remains "not applicable" (unchanged)

# --- Row 27 (Shooting Vic.) -------------------------------------------------------
# C27 stays "0" (text, unchanged)
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
# F27 goes from a number (1) to the literal text "0"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -100
# I27 stays "0" (text, unchanged)
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -100
$ws.Range("L27").Value = -100
# M27, N27 remain "not applicable" (unchanged)

# --- Row 28 (Shooting Inc.) — L28 unchanged numerically, N28 becomes numeric ---
$ws.Range("L28").Value = -100
$ws.Range("N28").NumberFormat = $fmtPct
$ws.Range("N28").Value = -100

# --- Row 29 (Hate Crimes) — L29 unchanged numerically, N29 becomes numeric -----
$ws.Range("L29").Value = -100
$ws.Range("N29").NumberFormat = $fmtPct
$ws.Range("N29").Value = -100
